$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet2 header row ----
$ws2.Range("A1").Value = "Age"
$ws2.Range("B1").Value = "Type 1 high round"
$ws2.Range("C1").Value = "Type 2 low round"
$ws2.Range("D1").Value = "Type 3  high peak"
$ws2.Range("E1").Value = "Type 4 low peak"

# ---- Row 3 : all zero ----
$ws2.Range("A3").Value = 0
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0

# ---- Row 4 : Age formula, rest zero ----
$ws2.Range("A4").Formula = "=A3+0.5"
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0

# ---- Age column, filled down as one shared formula block A5:A45 ----
$ws2.Range("A5:A45").Formula = "=A4+0.5"

# ---- Rows 5-14 : zero data ----
for ($r = 5; $r -le 14; $r++) {
    $ws2.Range("B$r").Value = 0
    $ws2.Range("C$r").Value = 0
    $ws2.Range("D$r").Value = 0
    $ws2.Range("E$r").Value = 0
}

# ---- Rows 15-43 : the actual life-table data ----
$bVals = @{
    15=30; 16=48; 17=61; 18=65; 19=68; 20=71; 21=73; 22=75; 23=76; 24=77;
    25=78; 26=79; 27=79.5; 28=80; 29=80; 30=80; 31=79.5; 32=79; 33=78; 34=77;
    35=76; 36=75; 37=73; 38=71; 39=68; 40=65; 41=61; 42=48; 43=30
}
$dVals = @{
    15=2; 16=4; 17=6; 18=9; 19=12; 20=16; 21=20; 22=25; 23=30; 24=40;
    25=50; 26=60; 27=75; 28=90; 29=96; 30=100; 31=96; 32=90; 33=75; 34=60;
    35=50; 36=40; 37=30; 38=22; 39=16; 40=11; 41=8; 42=4; 43=2
}

for ($r = 15; $r -le 43; $r++) {
    $ws2.Range("B$r").Value = $bVals[$r]
    $ws2.Range("D$r").Value = $dVals[$r]
}

# C column (=B/2) and E column (=D/2) - fill as shared formula blocks
$ws2.Range("C16:C43").Formula = "=B16/2"
$ws2.Range("E16:E43").Formula = "=D16/2"
# Row 15 gets its own (non-shared) formula, matching the anchor pattern
$ws2.Range("C15").Formula = "=B15/2"
$ws2.Range("E15").Formula = "=D15/2"
# Row 21 in column C is re-entered individually (breaks out of the shared run)
$ws2.Range("C21").Formula = "=B21/2"

# ---- Rows 44-45 : zero data, Age already filled above ----
for ($r = 44; $r -le 45; $r++) {
    $ws2.Range("B$r").Value = 0
    $ws2.Range("C$r").Value = 0
    $ws2.Range("D$r").Value = 0
    $ws2.Range("E$r").Value = 0
}

# ---- Select & activate Sheet2 (this becomes the active tab) ----
$ws2.Activate() | Out-Null
$ws2.Range("I12").Select() | Out-Null

Write-Output "edit complete"
